$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.3547408878803253
$ws.Range("B1").Value = 3.095179080963135
$ws.Range("C1").Value = 6.00098991394043
$ws.Range("D1").Value = 1.686342716217041
$ws.Range("E1").Value = 1.00642204284668
